# Add a new "DefaultHitTime" row to the Skill.xlsx "Property" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

$ws.Range("A$row").Value = "DefaultHitTime"
$ws.Range("B$row").Value = "float"
$ws.Range("C$row").Value = $false
$ws.Range("D$row").Value = $false
$ws.Range("E$row").Value = $false
$ws.Range("F$row").Value = $true
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = 0
$ws.Range("I$row").Value = "Friend"
$ws.Range("J$row").Value = "缺省打击时间（本来应该打到但是物理没碰撞到或者其他原因）"

# Match the text-formatted style used by the rest of column A/B/I/J (numFmtId 49 "@").
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("I$row").NumberFormat = "@"
$ws.Range("J$row").NumberFormat = "@"

# Keep the active-cell selection / top-left the same as the final workbook state.
$ws.Range("J32").Select() | Out-Null
